$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: set a cell's value as LITERAL TEXT, even if it looks like a date
# (e.g. "2024-11-08"), avoiding Excel's automatic date-serial coercion.
# We briefly force a text number format so the value is stored verbatim,
# then clear the format again so the cell ends up with the default style
# (matching the rest of the sheet, which carries no explicit style on the
# data columns).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ===========================================================================
# Sheet "展览" — plain numeric updates, no structural changes
# ===========================================================================
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8067
$ws1.Range("F5").Value = 5882
$ws1.Range("F6").Value = 500
$ws1.Range("F8").Value = 15
$ws1.Range("F9").Value = 70
$ws1.Range("F11").Value = 413

# ===========================================================================
# Sheet "演出" — a new event ("豫章D乐团") was added on 2024-11-08, sorted
# in just before the existing 2024-11-09 event, so a new row is inserted at
# row 4 and everything from the old row 4 onward shifts down by one.
# ===========================================================================
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(4).Insert()

# The blank row Insert() creates doesn't reliably pick up the "index column"
# border/alignment formatting used throughout column A, so copy it explicitly
# from the row above before filling in values.
$ws2.Range("A3:I3").Copy()
$ws2.Range("A4:I4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 4: the newly added event
$ws2.Range("A4").Value = 3
Set-TextValue $ws2.Range("B4") "2024-11-08"
$ws2.Range("C4").Value = "合肥·豫章D乐团-《蓬莱乐，万物生》——传统×先锋 疗愈音乐会"
$ws2.Range("D4").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws2.Range("E4").Value = "2024.11.08 19:30-11.08 21:00"
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 86.09999999999999
$ws2.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=92957"
$ws2.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202409/uifvAByr1727253170481.jpeg"

# Row 5 now holds what used to be row 4 (renumbered A5 = 4); all other
# columns are untouched from the pre-edit row 4.
$ws2.Range("A5").Value = 4

# Row 6 now holds what used to be row 5 (renumbered A6 = 5); all other
# columns are untouched from the pre-edit row 5.
$ws2.Range("A6").Value = 5

# ===========================================================================
# Sheet "全部类型" — combined view: same numeric updates as "展览" for the
# shared rows, plus the same new-row insertion as "演出" for the shared rows.
# ===========================================================================
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8067
$ws4.Range("F5").Value = 5882
$ws4.Range("F6").Value = 500
$ws4.Range("F8").Value = 15
$ws4.Range("F9").Value = 70

$ws4.Rows.Item(13).Insert()

# Same formatting fix as above for the "演出" sheet's inserted row.
$ws4.Range("A12:I12").Copy()
$ws4.Range("A13:I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 13: the newly added event
$ws4.Range("A13").Value = 12
Set-TextValue $ws4.Range("B13") "2024-11-08"
$ws4.Range("C13").Value = "合肥·豫章D乐团-《蓬莱乐，万物生》——传统×先锋 疗愈音乐会"
$ws4.Range("D13").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws4.Range("E13").Value = "2024.11.08 19:30-11.08 21:00"
$ws4.Range("F13").Value = 0
$ws4.Range("G13").Value = 86.09999999999999
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=92957"
$ws4.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202409/uifvAByr1727253170481.jpeg"

# Row 14 now holds what used to be row 13 (renumbered A14 = 13).
$ws4.Range("A14").Value = 13

# Row 15 now holds what used to be row 14 (renumbered A15 = 14), and its
# "想去人数" (F) value is updated the same way as the "展览" sheet's F11.
$ws4.Range("A15").Value = 14
$ws4.Range("F15").Value = 413

# Row 16 now holds what used to be row 15 (renumbered A16 = 15).
$ws4.Range("A16").Value = 15

# Row 17 now holds what used to be row 16 (renumbered A17 = 16).
$ws4.Range("A17").Value = 16
